# Scheduled market-data refresh: update currentAveragePrice / LevePrice /
# LeveProfit columns (H:N) for the affected leve rows on each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 863.5
$ws.Range("I28").Value = 836.2
$ws.Range("K28").Value = 836.2
$ws.Range("M28").Value = -351.2

$ws.Range("H33").Value = 1241.8889
$ws.Range("I33").Value = 993.25
$ws.Range("K33").Value = 993.25
$ws.Range("M33").Value = -764.25

$ws.Range("H62").Value = 4687.143
$ws.Range("I62").Value = 3782.2222
$ws.Range("K62").Value = 3782.2222
$ws.Range("M62").Value = -3158.2222

$ws.Range("H65").Value = 4687.143
$ws.Range("I65").Value = 3782.2222
$ws.Range("K65").Value = 18911.111
$ws.Range("M65").Value = -15791.111

$ws.Range("H103").Value = 707.8125
$ws.Range("I103").Value = 627.44446
$ws.Range("J103").Value = 811.1429000000001
$ws.Range("K103").Value = 1882.33338
$ws.Range("L103").Value = 2433.4287
$ws.Range("M103").Value = -1296.33338
$ws.Range("N103").Value = -3605.4287

$ws.Range("H132").Value = 3490.4375
$ws.Range("I132").Value = 3490.4375
$ws.Range("K132").Value = 10471.3125
$ws.Range("M132").Value = -7941.3125

$ws.Range("H137").Value = 8552.25
$ws.Range("I137").Value = 2721.1
$ws.Range("J137").Value = 11791.777
$ws.Range("K137").Value = 8163.299999999999
$ws.Range("L137").Value = 35375.331
$ws.Range("M137").Value = -5613.299999999999
$ws.Range("N137").Value = -40475.331

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7941389.5
$ws.Range("I32").Value = 8476347
$ws.Range("K32").Value = 8476347
$ws.Range("M32").Value = -8476060

$ws.Range("H46").Value = 33166.668
$ws.Range("J46").Value = 37750
$ws.Range("L46").Value = 37750
$ws.Range("N46").Value = -38388

$ws.Range("H61").Value = 38549660
$ws.Range("I61").Value = 125008504
$ws.Range("K61").Value = 125008504
$ws.Range("M61").Value = -125008292

$ws.Range("H74").Value = 7148902
$ws.Range("I74").Value = 13159152
$ws.Range("K74").Value = 13159152
$ws.Range("M74").Value = -13158278

$ws.Range("H77").Value = 7148902
$ws.Range("I77").Value = 13159152
$ws.Range("K77").Value = 65795760
$ws.Range("M77").Value = -65791392

$ws.Range("H108").Value = 44179.6
$ws.Range("J108").Value = 39947.5
$ws.Range("L108").Value = 39947.5
$ws.Range("N108").Value = -47627.5

$ws.Range("H136").Value = 38549660
$ws.Range("I136").Value = 125008504
$ws.Range("K136").Value = 375025512
$ws.Range("M136").Value = -375022962

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1602
$ws.Range("I105").Value = 527.5
$ws.Range("K105").Value = 527.5
$ws.Range("M105").Value = 1219.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1383314.8
$ws.Range("I31").Value = 33300.4
$ws.Range("K31").Value = 33300.4
$ws.Range("M31").Value = -33005.4

$ws.Range("H34").Value = 1383314.8
$ws.Range("I34").Value = 33300.4
$ws.Range("K34").Value = 33300.4
$ws.Range("M34").Value = -33098.4

$ws.Range("H132").Value = 2995.5293
$ws.Range("I132").Value = 2701.6333
$ws.Range("K132").Value = 8104.8999
$ws.Range("M132").Value = -5574.8999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 604.2222
$ws.Range("I14").Value = 604.2222
$ws.Range("K14").Value = 1812.6666
$ws.Range("M14").Value = -1639.6666

$ws.Range("H22").Value = 1849.5
$ws.Range("I22").Value = 1849.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 5548.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -5379.5
$ws.Range("N22").ClearContents()

$ws.Range("H24").Value = 8937.5
$ws.Range("I24").Value = 875
$ws.Range("K24").Value = 2625
$ws.Range("M24").Value = -2395

$ws.Range("H27").Value = 1849.5
$ws.Range("I27").Value = 1849.5
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 5548.5
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -5446.5
$ws.Range("N27").ClearContents()

$ws.Range("H107").Value = 562.36365
$ws.Range("I107").Value = 562.36365
$ws.Range("K107").Value = 1687.09095
$ws.Range("M107").Value = 232.90905

$ws.Range("H113").Value = 1300.7307
$ws.Range("I113").Value = 691.7778
$ws.Range("J113").Value = 1623.1177
$ws.Range("K113").Value = 2075.3334
$ws.Range("L113").Value = 4869.3531
$ws.Range("M113").Value = 94.66660000000002
$ws.Range("N113").Value = -9209.3531

$ws.Range("H137").Value = 6242.4443
$ws.Range("I137").Value = 6151.4
$ws.Range("J137").Value = 6356.25
$ws.Range("K137").Value = 18454.2
$ws.Range("L137").Value = 19068.75
$ws.Range("M137").Value = -13354.2
$ws.Range("N137").Value = -29268.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 93310.164
$ws.Range("J69").Value = 93310.164
$ws.Range("L69").Value = 93310.164
$ws.Range("N69").Value = -94808.164

$ws.Range("H70").Value = 4080
$ws.Range("I70").Value = 3850
$ws.Range("K70").Value = 3850
$ws.Range("M70").Value = -3580

$ws.Range("H72").Value = 93310.164
$ws.Range("J72").Value = 93310.164
$ws.Range("L72").Value = 279930.492
$ws.Range("N72").Value = -287418.492

$ws.Range("H73").Value = 4080
$ws.Range("I73").Value = 3850
$ws.Range("K73").Value = 3850
$ws.Range("M73").Value = -2914

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 40001876
$ws.Range("I93").Value = 52632900
$ws.Range("K93").Value = 52632900
$ws.Range("M93").Value = -52631652

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 10000
$ws.Range("I10").Value = 10000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 10000
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -9831
$ws.Range("N10").ClearContents()

$ws.Range("H113").Value = 975.5909
$ws.Range("I113").Value = 478.14285
$ws.Range("J113").Value = 1846.125
$ws.Range("K113").Value = 1434.42855
$ws.Range("L113").Value = 5538.375
$ws.Range("M113").Value = 735.5714499999999
$ws.Range("N113").Value = -9878.375

$ws.Range("H126").Value = 6117.95
$ws.Range("I126").Value = 5136.5713
$ws.Range("J126").Value = 8407.833000000001
$ws.Range("K126").Value = 15409.7139
$ws.Range("L126").Value = 25223.499
$ws.Range("M126").Value = -12939.7139
$ws.Range("N126").Value = -30163.499

$ws.Range("H136").Value = 2942.2693
$ws.Range("I136").Value = 2002.7646
$ws.Range("K136").Value = 6008.293799999999
$ws.Range("M136").Value = -3458.293799999999
